$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sprint "2017.10.06" sheet (sheet12): fill in the "Worked" (E) column for
#    the "Design responsive" row, which recalculates the Remain column and
#    the summary totals.
# ---------------------------------------------------------------------------
$ws1006 = $wb.Worksheets.Item("2017.10.06")
$ws1006.Range("E3").Value = 6
$ws1006.Range("G5").Select()

# ---------------------------------------------------------------------------
# 2. Add the new sprint sheet "2017.17.10" right after "2017.10.06".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1006)
$newSheet.Name = "2017.17.10"
$newSheet.PageSetup.TopMargin = 56.692913399999995
$newSheet.PageSetup.BottomMargin = 56.692913399999995

# Header row (reuses the existing shared strings from the other sprint sheets)
$newSheet.Range("A1").Value = "User Story"
$newSheet.Range("B1").Value = "Task"
$newSheet.Range("C1").Value = "Initial Estimation"
$newSheet.Range("D1").Value = "Current Estimation"
$newSheet.Range("E1").Value = "Worked"
$newSheet.Range("F1").Value = "Remain"

# User story block: "All" / "Model Change", "Statistics", "Tutor accept"
$newSheet.Range("A2").Value = "All"

$newSheet.Range("B3").Value = "Model Change"
$newSheet.Range("C3").Value = 1
$newSheet.Range("D3").Value = 1
$newSheet.Range("F3").Formula = "=D3-E3"

$newSheet.Range("B4").Value = "Statistics"
$newSheet.Range("C4").Value = 2
$newSheet.Range("D4").Value = 2

$newSheet.Range("B5").Value = "Tutor accept"
$newSheet.Range("C5").Value = 4
$newSheet.Range("D5").Value = 4

# F4 and F5 share one formula group (mirrors a fill-down in Excel).
$newSheet.Range("F4:F5").Formula = "=D4-E4"

# Second user story block: "All" / "Design responsive"
$newSheet.Range("A8").Value = "All"
$newSheet.Range("B8").Value = "Design responsive"
$newSheet.Range("C8").Value = 6
$newSheet.Range("D8").Value = 6
$newSheet.Range("F8").Formula = "=D8-E8"

# Totals row
$newSheet.Range("C12").Formula = "=SUM(C2:C11)"
$newSheet.Range("D12").Formula = "=SUM(D2:D11)"
$newSheet.Range("E12").Formula = "=SUM(E2:E11)"
$newSheet.Range("F12").Formula = "=SUM(F2:F11)"

# Planned hours block
$newSheet.Range("A14").Value = "Name"
$newSheet.Range("B14").Value = "Planned hours"
$newSheet.Range("A15").Value = "Eva"
$newSheet.Range("B15").Value = 7
$newSheet.Range("A16").Value = "Danijal"
$newSheet.Range("B16").Value = 6
$newSheet.Range("B17").Formula = "=SUM(B15:B16)"

$newSheet.Range("B15").Select()

# ---------------------------------------------------------------------------
# 3. Summary sheet: log the new sprint and make it the active/selected tab.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A16").Value = 43025
$summary.Range("A15").Copy()
$summary.Range("A16").PasteSpecial(-4122)
$summary.Application.CutCopyMode = $false
$summary.Range("B16").Value = "Statistics"
$summary.Range("B16").Select()
$summary.Activate()
